$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the shared string "Equitable Service" -> "EquiBliss" (header in J2)
$ws.Range("J2").Value = "EquiBliss"

# Update the raw scheduler data in column I (rows 3-6); dependent formulas
# in I8, I9, I10, I12, I13, I14 recalc automatically.
$ws.Range("I3").Value = 1609951
$ws.Range("I4").Value = 25259910
$ws.Range("I5").Value = 23822041
$ws.Range("I6").Value = 4637496

# Update the active selection shown in the sheet view.
$ws.Range("L7").Select()
